$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.430374
$ws.Range("H2").Value = 7.291122000000001
$ws.Range("I2").Value = 0.009222757332915244
$ws.Range("J2").Value = 0.009222757332915246
$ws.Range("M2").Value = 1.533541666666667
$ws.Range("N2").Value = 4.600625
$ws.Range("O2").Value = 0.01998214594581092
$ws.Range("P2").Value = 0.01998214594581093
$ws.Range("Q2").Value = 3.727079794583333
$ws.Range("R2").Value = 33.54371815125
$ws.Range("S2").Value = 0.0001842904830491103
$ws.Range("T2").Value = 0.0001842904830491104
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.430374
$ws.Range("H3").Value = 7.291122000000001
$ws.Range("I3").Value = 0.009222757332915244
$ws.Range("J3").Value = 0.009222757332915246
$ws.Range("M3").Value = 3.948587333333334
$ws.Range("O3").Value = 0.05145034536032411
$ws.Range("P3").Value = 0.05145034536032412
$ws.Range("Q3").Value = 9.596543991662667
$ws.Range("R3").Value = 86.36889592496401
$ws.Range("S3").Value = 0.000474514049952951
$ws.Range("T3").Value = 0.0004745140499529512
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.430374
$ws.Range("H4").Value = 7.291122000000001
$ws.Range("I4").Value = 0.009222757332915244
$ws.Range("J4").Value = 0.009222757332915246
$ws.Range("M4").Value = 70.69501233333334
$ws.Range("N4").Value = 212.085037
$ws.Range("O4").Value = 0.921160529766436
$ws.Range("P4").Value = 0.9211605297664361
$ws.Range("Q4").Value = 171.8153199046127
$ws.Range("R4").Value = 1546.337879141514
$ws.Range("S4").Value = 0.008495640030695488
$ws.Range("T4").Value = 0.008495640030695492
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.430374
$ws.Range("H5").Value = 7.291122000000001
$ws.Range("I5").Value = 0.009222757332915244
$ws.Range("J5").Value = 0.009222757332915246
$ws.Range("M5").Value = 0.568453
$ws.Range("N5").Value = 1.705359
$ws.Range("O5").Value = 0.007406978927428811
$ws.Range("P5").Value = 0.007406978927428812
$ws.Range("Q5").Value = 1.381553391422
$ws.Range("R5").Value = 12.433980522798
$ws.Range("S5").Value = [double]"6.831276921769275E-05"
$ws.Range("T5").Value = [double]"6.831276921769278E-05"
$ws.Range("I6").Value = 0.5480399755605952
$ws.Range("J6").Value = 0.5480399755605954
$ws.Range("M6").Value = 1.533541666666667
$ws.Range("N6").Value = 4.600625
$ws.Range("O6").Value = 0.01998214594581092
$ws.Range("P6").Value = 0.01998214594581093
$ws.Range("Q6").Value = 221.4726730634027
$ws.Range("R6").Value = 1993.254057570625
$ws.Range("S6").Value = 0.01095101477579047
$ws.Range("T6").Value = 0.01095101477579047
$ws.Range("I7").Value = 0.5480399755605952
$ws.Range("J7").Value = 0.5480399755605954
$ws.Range("M7").Value = 3.948587333333334
$ws.Range("O7").Value = 0.05145034536032411
$ws.Range("P7").Value = 0.05145034536032412
$ws.Range("S7").Value = 0.02819684601385621
$ws.Range("T7").Value = 0.02819684601385622
$ws.Range("I8").Value = 0.5480399755605952
$ws.Range("J8").Value = 0.5480399755605954
$ws.Range("M8").Value = 70.69501233333334
$ws.Range("N8").Value = 212.085037
$ws.Range("O8").Value = 0.921160529766436
$ws.Range("P8").Value = 0.9211605297664361
$ws.Range("Q8").Value = 10209.70847681362
$ws.Range("R8").Value = 91887.37629132262
$ws.Range("S8").Value = 0.5048327942205826
$ws.Range("T8").Value = 0.5048327942205827
$ws.Range("I9").Value = 0.5480399755605952
$ws.Range("J9").Value = 0.5480399755605954
$ws.Range("M9").Value = 0.568453
$ws.Range("N9").Value = 1.705359
$ws.Range("O9").Value = 0.007406978927428811
$ws.Range("P9").Value = 0.007406978927428812
$ws.Range("Q9").Value = 82.09545795684966
$ws.Range("R9").Value = 738.859121611647
$ws.Range("S9").Value = 0.004059320550365929
$ws.Range("T9").Value = 0.004059320550365931
$ws.Range("G10").Value = 116.470388
$ws.Range("H10").Value = 349.411164
$ws.Range("I10").Value = 0.4419805861132828
$ws.Range("J10").Value = 0.4419805861132828
$ws.Range("M10").Value = 1.533541666666667
$ws.Range("N10").Value = 4.600625
$ws.Range("O10").Value = 0.01998214594581092
$ws.Range("P10").Value = 0.01998214594581093
$ws.Range("Q10").Value = 178.6121929308333
$ws.Range("R10").Value = 1607.5097363775
$ws.Range("S10").Value = 0.008831720576930669
$ws.Range("T10").Value = 0.008831720576930671
$ws.Range("G11").Value = 116.470388
$ws.Range("H11").Value = 349.411164
$ws.Range("I11").Value = 0.4419805861132828
$ws.Range("J11").Value = 0.4419805861132828
$ws.Range("M11").Value = 3.948587333333334
$ws.Range("O11").Value = 0.05145034536032411
$ws.Range("P11").Value = 0.05145034536032412
$ws.Range("Q11").Value = 459.8934987652187
$ws.Range("R11").Value = 4139.041488886968
$ws.Range("S11").Value = 0.02274005379808687
$ws.Range("T11").Value = 0.02274005379808688
$ws.Range("G12").Value = 116.470388
$ws.Range("H12").Value = 349.411164
$ws.Range("I12").Value = 0.4419805861132828
$ws.Range("J12").Value = 0.4419805861132828
$ws.Range("M12").Value = 70.69501233333334
$ws.Range("N12").Value = 212.085037
$ws.Range("O12").Value = 0.921160529766436
$ws.Range("P12").Value = 0.9211605297664361
$ws.Range("Q12").Value = 8233.875516128119
$ws.Range("R12").Value = 74104.87964515307
$ws.Range("S12").Value = 0.4071350708505915
$ws.Range("T12").Value = 0.4071350708505916
$ws.Range("G13").Value = 116.470388
$ws.Range("H13").Value = 349.411164
$ws.Range("I13").Value = 0.4419805861132828
$ws.Range("J13").Value = 0.4419805861132828
$ws.Range("M13").Value = 0.568453
$ws.Range("N13").Value = 1.705359
$ws.Range("O13").Value = 0.007406978927428811
$ws.Range("P13").Value = 0.007406978927428812
$ws.Range("Q13").Value = 66.20794146976399
$ws.Range("R13").Value = 595.871473227876
$ws.Range("S13").Value = 0.003273740887673721
$ws.Range("T13").Value = 0.003273740887673721
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.1994
$ws.Range("H14").Value = 0.5982000000000001
$ws.Range("I14").Value = 0.0007566809932065188
$ws.Range("J14").Value = 0.0007566809932065189
$ws.Range("M14").Value = 1.533541666666667
$ws.Range("N14").Value = 4.600625
$ws.Range("O14").Value = 0.01998214594581092
$ws.Range("P14").Value = 0.01998214594581093
$ws.Range("Q14").Value = 0.3057882083333334
$ws.Range("R14").Value = 2.752093875
$ws.Range("S14").Value = [double]"1.512011004067382E-05"
$ws.Range("T14").Value = [double]"1.512011004067383E-05"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.1994
$ws.Range("H15").Value = 0.5982000000000001
$ws.Range("I15").Value = 0.0007566809932065188
$ws.Range("J15").Value = 0.0007566809932065189
$ws.Range("M15").Value = 3.948587333333334
$ws.Range("O15").Value = 0.05145034536032411
$ws.Range("P15").Value = 0.05145034536032412
$ws.Range("Q15").Value = 0.7873483142666668
$ws.Range("R15").Value = 7.086134828400001
$ws.Range("S15").Value = [double]"3.893149842806846E-05"
$ws.Range("T15").Value = [double]"3.893149842806847E-05"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.1994
$ws.Range("H16").Value = 0.5982000000000001
$ws.Range("I16").Value = 0.0007566809932065188
$ws.Range("J16").Value = 0.0007566809932065189
$ws.Range("M16").Value = 70.69501233333334
$ws.Range("N16").Value = 212.085037
$ws.Range("O16").Value = 0.921160529766436
$ws.Range("P16").Value = 0.9211605297664361
$ws.Range("Q16").Value = 14.09658545926667
$ws.Range("R16").Value = 126.8692691334
$ws.Range("S16").Value = 0.0006970246645663098
$ws.Range("T16").Value = 0.00069702466456631
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.1994
$ws.Range("H17").Value = 0.5982000000000001
$ws.Range("I17").Value = 0.0007566809932065188
$ws.Range("J17").Value = 0.0007566809932065189
$ws.Range("M17").Value = 0.568453
$ws.Range("N17").Value = 1.705359
$ws.Range("O17").Value = 0.007406978927428811
$ws.Range("P17").Value = 0.007406978927428812
$ws.Range("Q17").Value = 0.1133495282
$ws.Range("R17").Value = 1.0201457538
$ws.Range("S17").Value = [double]"5.604720171466588E-06"
$ws.Range("T17").Value = [double]"5.604720171466589E-06"
